## Applies the edit described by the diff:
##  - Inserts a new slide (Title + Content layout) at position 5, pushing the
##    existing "Bye world" slide (was slide 5) down to position 6.
##  - The new slide stays blank (title/content placeholders empty), same as
##    the captured target (sldId 261).
##  - The new slide's notes page gets the "BREAK TITLEEE" placeholder text
##    that was used to author a break-title slide.
##  - The "Bye world" slide (now slide 6) keeps its original notes text
##    ("Skipped previous comment") untouched.

$p = $ppt.ActivePresentation

# Insert a new "Title and Content" slide right before the current slide 5
# ("Bye world"), so it becomes the new slide 5 and "Bye world" becomes 6.
$newSlide = $p.Slides.Add(5, 2)

# Give the new slide's notes page the break-title JSON text.
$notes = $newSlide.NotesPage
$notesBody = $notes.Shapes.Item(2)
$notesBody.TextFrame.TextRange.Text = "{" + [char]0x201C + "display_name" + [char]0x201D + ": " + [char]0x201C + "BREAK TITLEEE" + [char]0x201D + "}"
